$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.332.30"
$ws.Range("D3").Value = "'1.760.34"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'303.94"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("D7").Value = "'0.4276"
$ws.Range("D8").Value = "'0.3615"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("D9").Value = "'0.07052"
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").Value = "'0.8291"
$ws.Range("E10").Value = "  -4.11%  "
$ws.Range("D11").Value = "'20.12"
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("D12").Value = "'1.750.70"
$ws.Range("E12").Value = "  -4.11%  "
$ws.Range("D13").Value = "'5.226"
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("D14").Value = "'6.363"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "'0.06798"
$ws.Range("E15").Value = "  -2.18%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'79.02"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("D18").Value = "'0.000008649"
$ws.Range("E18").Value = "  -2.88%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'14.91"
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("D21").Value = "'26.129.20"
$ws.Range("E21").Value = "  -5.22%  "
$ws.Range("D22").Value = "'4.984"
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("D23").Value = "'11.07"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").Value = "'1.965.20"
$ws.Range("E24").Value = "  -4.98%  "
$ws.Range("D25").Value = "'1.904"
$ws.Range("E25").Value = "  -4.48%  "
$ws.Range("D26").Value = "'151.82"
$ws.Range("E26").Value = "  -1.85%  "
$ws.Range("D27").Value = "'18.09"
$ws.Range("E27").Value = "  -4.28%  "
$ws.Range("D28").Value = "'5.024"
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("D29").Value = "'114.10"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "'1.670"
$ws.Range("E30").Value = "  -8.42%  "
$ws.Range("D31").Value = "'0.08867"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").Value = "'0.7221"
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.104"
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.298"
$ws.Range("E34").Value = "  -5.40%  "
$ws.Range("D35").Value = "'1.000"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'2.703"
$ws.Range("E36").Value = "  -9.86%  "
$ws.Range("D37").Value = "'1.067"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("D38").Value = "'0.05092"
$ws.Range("E38").Value = "  -4.41%  "
$ws.Range("D39").Value = "'0.01876"
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("D40").Value = "'0.1601"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").Value = "'0.4884"
$ws.Range("E41").Value = "  -3.88%  "
$ws.Range("D42").Value = "'2.485"
$ws.Range("E42").Value = "  -11.22%  "
$ws.Range("D43").Value = "'6.139"
$ws.Range("E43").Value = "  -5.37%  "
$ws.Range("D44").Value = "'7.981"
$ws.Range("E44").Value = "  -4.27%  "
$ws.Range("D45").Value = "'104.37"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "'9.972"
$ws.Range("E47").Value = "  -3.93%  "
$ws.Range("D48").Value = "'0.06180"
$ws.Range("E48").Value = "  -4.51%  "
$ws.Range("D49").Value = "'0.4461"
$ws.Range("E49").Value = "  -4.53%  "
$ws.Range("D50").Value = "'1.564"
$ws.Range("E50").Value = "  -3.38%  "
$ws.Range("D51").Value = "'1.708"
$ws.Range("E51").Value = "  -1.57%  "
